# Added creation of table for primer: fill in Thermocycler, Freeze and
# Squeeze, Specimen Discs, Kim Wipe, and Thiol-Streptavidin rows with
# company/product-number data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: Thermocycler -> Bio-Rad / T100
$ws.Range("B31").Value = "Bio-Rad"
$ws.Range("C31").Value = "T100"

# Row 34: Kim Wipe -> KimTech / 5511
$ws.Range("B34").Value = "KimTech"
$ws.Range("C34").Value = 5511

# Row 32: Freeze and Squeeze -> Bio-Rad / 732-6165
$ws.Range("B32").Value = "Bio-Rad"
$ws.Range("C32").Value = "732-6165"

# Row 25: Thiol-Streptavidin -> Protein Mods / SAVT
$ws.Range("B25").Value = "Protein Mods"
$ws.Range("C25").Value = "SAVT"

# Row 33: Specimen Discs, Ted Pella -> 16218
$ws.Range("C33").Value = 16218

# Move the active selection to E29, matching the last user interaction.
$ws.Range("E29").Select() | Out-Null
